$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 138, shifting existing rows 138-156 down to 139-157
$ws.Rows("138:138").Insert()

# Populate the newly inserted row 138 with the new data record
$ws.Cells.Item(138, 1).Value = 10
$ws.Cells.Item(138, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(138, 3).Value = "La Araucanía"
$ws.Cells.Item(138, 4).Value = 44491
$ws.Cells.Item(138, 4).NumberFormat = $ws.Cells.Item(137, 4).NumberFormat
$ws.Cells.Item(138, 5).Value = 9
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100102
$ws.Cells.Item(138, 8).Value = "Cítricos"
$ws.Cells.Item(138, 9).Value = 100102006
$ws.Cells.Item(138, 10).Value = "Pomelo"
$ws.Cells.Item(138, 11).Value = "Start Ruby"
$ws.Cells.Item(138, 12).Value = "Primera"
$ws.Cells.Item(138, 13).Value = 90
$ws.Cells.Item(138, 14).Value = 10000
$ws.Cells.Item(138, 15).Value = 10000
$ws.Cells.Item(138, 16).Value = 10000
$ws.Cells.Item(138, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(138, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(138, 19).Value = 667
$ws.Cells.Item(138, 20).Value = 15

Write-Output "Row inserted and populated"
